# Reorders the header columns and updates the one-hot selection flags in
# 90_scenecat_block_order.xlsx (sheet1) to the new block-order layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header labels (columns re-ordered) ---
$ws.Range("A1").Value = "kitchens_1"
$ws.Range("B1").Value = "bedrooms_1"
$ws.Range("C1").Value = "living_rooms_1"
$ws.Range("D1").Value = "living_rooms_2"
$ws.Range("E1").Value = "kitchens_2"
$ws.Range("F1").Value = "bedrooms_2"

# --- Rows 2-7: one-hot indicator values per row ---
$data = @(
    @(0,1,0,0,0,0),
    @(0,0,0,0,1,0),
    @(0,0,0,1,0,0),
    @(0,0,0,0,0,1),
    @(1,0,0,0,0,0),
    @(0,0,1,0,0,0)
)

$cols = @("A","B","C","D","E","F")

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $cellRef = $cols[$j] + $rowNum
        $ws.Range($cellRef).Value = $rowVals[$j]
    }
}
